# Weekly update: add a new price-report date (2022-04-29, serial 44680) for
# "Pepino dulce" at Mercado Mayorista Lo Valledor de Santiago, inserted as
# three new rows right before the existing 2021-03-24 (serial 44279) block,
# pushing every subsequent row down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 144, shifting rows 144:230 down to 147:233.
$ws.Rows("144:146").Insert()

# Clone the row immediately below (now rows 147:149, which still hold the
# data formerly on 144:146) into the freshly inserted, empty rows so every
# column (A-R) is fully populated with the correct formatting/text, then
# only the columns that actually differ for the new date get overwritten.
$ws.Range("A147:R149").Copy() | Out-Null
$ws.Range("A144").PasteSpecial() | Out-Null

# Row 144: Especial
$ws.Range("D144").Value = 44680
$ws.Range("J144").Value = 210
$ws.Range("K144").Value = 18000
$ws.Range("L144").Value = 18000
$ws.Range("M144").Value = 18000
$ws.Range("P144").Value = 1000

# Row 145: Primera
$ws.Range("D145").Value = 44680
$ws.Range("J145").Value = 270
$ws.Range("K145").Value = 15000
$ws.Range("L145").Value = 15000
$ws.Range("M145").Value = 15000
$ws.Range("P145").Value = 833

# Row 146: Segunda
$ws.Range("D146").Value = 44680
$ws.Range("J146").Value = 220
$ws.Range("K146").Value = 13000
$ws.Range("L146").Value = 13000
$ws.Range("M146").Value = 13000
$ws.Range("P146").Value = 722
